$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 2 and add new rows 3-10 with the refreshed domain list ---
$ws.Range("A2").Value = "2025-1"
$ws.Range("B2").Value = "News-Media"
$ws.Range("C2").Value = "LV"
$ws.Range("D2").Value = "https://rebaltica.lv/"
$ws.Range("E2").Value = 45703

$ws.Range("A3").Value = "2025-2"
$ws.Range("B3").Value = "News-Media"
$ws.Range("C3").Value = "NEWS"
$ws.Range("D3").Value = "https://testpress.news/"
$ws.Range("E3").Value = 45703

$ws.Range("A4").Value = "2025-3"
$ws.Range("B4").Value = "Shopping"
$ws.Range("C4").Value = "UA"
$ws.Range("D4").Value = "https://kasta.ua/"
$ws.Range("E4").Value = 45703

$ws.Range("A5").Value = "2025-4"
$ws.Range("B5").Value = "Foto"
$ws.Range("C5").Value = "LV"
$ws.Range("D5").Value = "https://www.unfoto.lv/"
$ws.Range("E5").Value = 45703

$ws.Range("A6").Value = "2025-5"
$ws.Range("B6").Value = "Edu"
$ws.Range("C6").Value = "EDU"
$ws.Range("D6").Value = "https://www.harvard.edu/"
$ws.Range("E6").Value = 45703

$ws.Range("A7").Value = "2025-6"
$ws.Range("B7").Value = "Gov"
$ws.Range("C7").Value = "PL"
$ws.Range("D7").Value = "https://www.gov.pl/web/diplomacy"
$ws.Range("E7").Value = 45703

$ws.Range("A8").Value = "2025-7"
$ws.Range("B8").Value = "Adult"
$ws.Range("C8").Value = "COM"
$ws.Range("D8").Value = "https://fapello.com/"
$ws.Range("E8").Value = 45703

$ws.Range("A9").Value = "2025-8"
$ws.Range("B9").Value = "News-Media"
$ws.Range("C9").Value = "COM"
$ws.Range("D9").Value = "https://www.aljazeera.com/"
$ws.Range("E9").Value = 45705

$ws.Range("A10").Value = "2025-9"
$ws.Range("B10").Value = "News-Media"
$ws.Range("C10").Value = "COM"
$ws.Range("D10").Value = "https://www.bbc.com/"
$ws.Range("E10").Value = 45705

# --- Row heights shrink slightly across the whole used range ---
$ws.Range("A1:E16").EntireRow.RowHeight = 18.75

# --- Date column's alignment reverts from left to the default (general) ---
$ws.Range("E1:E16").HorizontalAlignment = 1
